$d = $word.ActiveDocument

# 1) Profile paragraph: "... Engineering student about to graduate from the
#    University of Essex ..." -> "... Engineering soon graduating from the
#    University of Essex ..." (typo/wording fix)
$d.Content.Find.Execute("student about to graduate", $true, $false, $false, $false, $false, $true, 1, $false, "soon graduating", 2) | Out-Null

# 2) Interests paragraph: break the run-on sentence into two lines. The text
#    stays the same ("... from across the world. I strive to solve real-world
#    problems ...") but a manual line break is inserted right after the space
#    that follows the period and before "I strive to solve real-world...".
$d.Content.Find.Execute(". I strive to solve real", $true, $false, $false, $false, $false, $true, 1, $false, ". ^lI strive to solve real", 2) | Out-Null
